$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.299146413803101
$ws.Range("B1").Value = 1.384466052055359
$ws.Range("C1").Value = 1.579964518547058
$ws.Range("D1").Value = 2.51872730255127
$ws.Range("E1").Value = 15
